$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.325.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "'2.253.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.26%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "'247.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.69%  "
$ws.Range("D6").Value = "'0.629"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.98%  "
$ws.Range("D7").Value = "'73.97"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.91%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.620"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.08%  "
$ws.Range("D10").Value = "'41.36"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.10%  "
$ws.Range("D11").Value = "'0.0941"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.01%  "
$ws.Range("D12").Value = "'7.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.68%  "
$ws.Range("E13").Value = "  -2.88%  "
$ws.Range("D14").Value = "'2.585.19"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").Value = "'14.50"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.40%  "
$ws.Range("D16").Value = "'0.854"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.07%  "
$ws.Range("D17").Value = "'2.245.63"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.23%  "
$ws.Range("D18").Value = "'42.102.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.0₃0980"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.48%  "
$ws.Range("D20").Value = "'6.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("D21").Value = "'71.87"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.82%  "
$ws.Range("E22").Value = "  +5.16%  "
$ws.Range("D23").Value = "'231.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("D25").Value = "'11.15"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("D26").Value = "'7.98"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +26.20%  "
$ws.Range("D27").Value = "'3.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.18%  "
$ws.Range("D28").Value = "'2.32"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.56%  "
$ws.Range("E29").Value = "  +2.70%  "
$ws.Range("D30").Value = "'169.02"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.71%  "
$ws.Range("D31").Value = "'20.69"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.70%  "
$ws.Range("E32").Value = "  -7.43%  "
$ws.Range("E33").Value = "  -5.97%  "
$ws.Range("D34").Value = "'30.30"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.64%  "
$ws.Range("E35").Value = "  -1.59%  "
$ws.Range("D36").Value = "'4.50"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.87%  "
$ws.Range("D37").Value = "'4.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.35%  "
$ws.Range("D38").Value = "'0.0301"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.54%  "
$ws.Range("D39").Value = "'13.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("E40").Value = "  -5.26%  "
$ws.Range("E41").Value = "  -2.04%  "
$ws.Range("D42").Value = "'61.99"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.18%  "
$ws.Range("D43").Value = "'0.203"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.44%  "
$ws.Range("D44").Value = "'108.44"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.04%  "
$ws.Range("E45").Value = "  -4.10%  "
$ws.Range("E46").Value = "  -0.24%  "
$ws.Range("D47").Value = "'0.996"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.49%  "
$ws.Range("E48").Value = "  -4.44%  "
$ws.Range("E49").Value = "  -1.27%  "
$ws.Range("E50").Value = "  +0.66%  "
$ws.Range("B51").Value = "HuobiToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D51").Value = "'2.70"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.80%  "
